# SP-29 BIS-28: add "Value Unit" property row to openbis-metadata sheet,
# fill in example Start Data Row/Col values, and drop the now-redundant
# "Value Unit" header row from openbis-data.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # openbis-metadata
$ws2 = $wb.Worksheets.Item(2)   # openbis-data

# --- sheet1 (openbis-metadata): insert a new "Value Unit" row above the
#     existing "Header Format" row (old row 6), pushing everything below
#     down by one. ---
$ws1.Rows.Item(6).Insert()

$ws1.Range("A6").Value = "Value Unit"
$ws1.Range("B6").Value = ""
$ws1.Range("C6").Value = "One of mM, uM, RatioT1, or RatioCs"
$ws1.Range("D6").Value = "mM"

# give C6 its own distinctive (italic, grey, 14pt) look + left border
$c = $ws1.Range("C6")
$c.Font.Name = "Verdana"
$c.Font.Italic = $true
$c.Font.Size = 14
$c.Font.Color = 8421504
$c.Borders.Item(7).LineStyle = 1

# fill in example values for the rows that shifted down
$ws1.Range("B7").Value = "METABOL HYBRID"
$ws1.Range("B8").Value = 3
$ws1.Range("D8").Value = 3
$ws1.Range("B9").Value = "C"

$ws1.Range("C13").Select()

# --- sheet2 (openbis-data): the "Value Unit" header row is no longer
#     needed (it now lives on sheet1), so remove it. ---
$ws2.Rows.Item(3).Delete()
$ws2.Range("A3:XFD3").Select()
